$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new row at position 7, shifting existing rows 7-21 down to 8-22
$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value = 18330051920213
$ws.Cells.Item(7, 2).Value = "LADINO"
$ws.Cells.Item(7, 3).Value = "URBINA"
$ws.Cells.Item(7, 4).Value = "MARIBEL"
$ws.Cells.Item(7, 5).Value = "TEMAS DE FILOSOFÍA"
$ws.Cells.Item(7, 6).Value = "6ARHM"
$ws.Cells.Item(7, 7).Value = 2
